$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.945.88"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "3.152.30"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.05"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.03"
$ws.Range("E6").Value = "  -1.09%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.144.35"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.22"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "3.672.80"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.30"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.148.23"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "63.843.82"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.81"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.39"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.733"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.01"
$ws.Range("E24").Value = "  -2.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  +5.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.32"
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.72"
$ws.Range("E28").Value = "  +7.88%  "
$ws.Range("E29").Value = "  +7.11%  "
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.78"
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").Value = "0.0₃0844"
$ws.Range("E35").Value = "  -4.72%  "
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("E37").Value = "  -3.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.17"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.27"
$ws.Range("E39").Value = "  -5.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "463.03"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.25"
$ws.Range("E42").Value = "  +5.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.294"
$ws.Range("E43").Value = "  +5.20%  "
$ws.Range("D44").Value = "2.928.31"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.13"
$ws.Range("E46").Value = "  +12.66%  "
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.66"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.26"
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("E51").Value = "  -1.31%  "
